$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Fix typo: "Group4094" -> "Groupp4094"
$ws.Range("B2").Value = "Groupp4094"

# Update the active selection to B2 (as reflected in the saved view state)
$ws.Activate()
$ws.Range("B2").Select()
